# Swap the order of "dnasr281@gmail.com" and "System" in the
# "Recorded By" column (G) of the Session Analysis Results sheet:
#   "dnasr281@gmail.com, System"  ->  "System, dnasr281@gmail.com"
#
# This text appears only in column G, so restrict the Find/Replace to
# that column to be safe and efficient.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = $ws.Columns.Item(7)  # Column G ("Recorded By")
$col.Replace("dnasr281@gmail.com, System", "System, dnasr281@gmail.com")
